$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.919.77"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.207.83"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'229.80"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").Value = "'60.63"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.401"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").Value = "'0.0900"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "2.538.35"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "'15.40"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").Value = "'21.99"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "'0.797"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "'5.56"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "2.230.62"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "41.877.98"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "0.0₃0938"
$ws.Range("E19").Value = "  +4.61%  "
$ws.Range("D20").Value = "'72.11"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "'6.08"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'242.54"
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "'2.39"
$ws.Range("E24").Value = "  +5.79%  "
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "'9.61"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").Value = "'168.88"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("D29").Value = "'20.22"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("E31").Value = "  -4.74%  "
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("E33").Value = "  -5.14%  "
$ws.Range("E34").Value = "  -3.24%  "
$ws.Range("D35").Value = "'0.0645"
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "'6.27"
$ws.Range("E36").Value = "  -6.65%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.53"
$ws.Range("E37").Value = "  -7.94%  "
$ws.Range("D38").Value = "'2.32"
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("B39").Value = "BinanceUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0243"
$ws.Range("E40").Value = "  +2.58%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.53"
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("B42").Value = "TerraClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D42").Value = "'0.000221"
$ws.Range("E42").Value = "  -11.87%  "
$ws.Range("E43").Value = "  -3.55%  "
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").Value = "'96.51"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("D46").Value = "'4.35"
$ws.Range("E46").Value = "  -13.56%  "
$ws.Range("D47").Value = "1.453.71"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").Value = "'16.05"
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("E50").Value = "  -3.56%  "
$ws.Range("E51").Value = "  +0.61%  "
